# correction de bug et redirection
# Populate the participants list for "La Descente des Alpes - M1" (2021-02-27)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$epreuve = "La Descente des Alpes - M1"
$date = "2021-02-27"
$type = "M1"

$data = @(
    @(1, "Génique", "Yoann"),
    @(2, "Mairot", "Jean-christophe"),
    @(3, "Cherief", "Saufiane"),
    @(4, "Rameau", "Célia"),
    @(5, "Ligourel", "Teedji"),
    @(6, "toto", "tata"),
    @(7, "choula", "poula"),
    @(8, "mairot", "tutu")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = 2
    $ws.Cells.Item($row, 5).Value = $epreuve

    # Enter the date with a leading apostrophe so it is kept as plain text
    # (otherwise it gets auto-converted to a date serial number), then clear
    # the resulting cell formatting so no extra style gets attached.
    $ws.Cells.Item($row, 6).Value = "'" + $date
    $ws.Cells.Item($row, 6).ClearFormats()

    $ws.Cells.Item($row, 10).Value = 1
    $ws.Cells.Item($row, 11).Value = $type

    $row = $row + 1
}
